$d = $word.ActiveDocument

function Split-RunAt($range) {
    # Force Word to keep `range` as its own run (not merged into a
    # neighboring run with identical formatting) by toggling a
    # character property on and back off.
    $range.Font.Bold = 1
    $range.Font.Bold = 0
}

# ---------------------------------------------------------------
# Edit 1: "The closest model was run 2/run 4 with either 100 or 200
# epochs; ". It did not achieve the desired accuracy of 75% but got
# closest." paragraph.
# ---------------------------------------------------------------

$full = $d.Content.Text
$pidx = $full.IndexOf("The closest model was run 2")

$prefix = "The closest model was run 2"
$seg1 = "/run 4 with either 100 or 200 epochs"
$seg2 = "; "
$seg3 = "i"
$seg4 = "t did not achieve the desired accuracy of 75% but got closest."

$runStart = $pidx + $prefix.Length
$paraEndMarker = "got closest."
$localEnd = $full.IndexOf($paraEndMarker, $pidx) + $paraEndMarker.Length

$r = $d.Range($runStart, $localEnd)
$r.Text = $seg1 + $seg2 + $seg3 + $seg4

$b0 = $runStart
$b1 = $b0 + $seg1.Length
$b2 = $b1 + $seg2.Length
$b3 = $b2 + $seg3.Length
$b4 = $b3 + $seg4.Length

Split-RunAt $d.Range($b0, $b1)
Split-RunAt $d.Range($b1, $b2)
Split-RunAt $d.Range($b2, $b3)
Split-RunAt $d.Range($b3, $b4)

# ---------------------------------------------------------------
# Edit 2: "The results of each of the model attempts ... Also, would
# recommend looking at changing the number of features." paragraph.
# ---------------------------------------------------------------

$full = $d.Content.Text
$pidx2 = $full.IndexOf("The results of each of the model attempts")

$seg1b = "The results of each of the model attempts did not improve the accuracy enough to reach the desired accuracy rate of 75%. In order to improve the model, it would be recommended to check for and remove outliers. Another consideration, that may "
$seg2b = "affect"
$seg3b = " results and accuracy, would be looking at changing the number of features."

$paraEndMarker2 = "changing the number of features."
$localEnd2 = $full.IndexOf($paraEndMarker2, $pidx2) + $paraEndMarker2.Length

$r2 = $d.Range($pidx2, $localEnd2)
$r2.Text = $seg1b + $seg2b + $seg3b

$c0 = $pidx2
$c1 = $c0 + $seg1b.Length
$c2 = $c1 + $seg2b.Length
$c3 = $c2 + $seg3b.Length

Split-RunAt $d.Range($c0, $c1)
Split-RunAt $d.Range($c1, $c2)
Split-RunAt $d.Range($c2, $c3)
